$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1312.3636
$ws.Range("I43").Value = 879
$ws.Range("K43").Value = 879
$ws.Range("M43").Value = -810
$ws.Range("H80").Value = 3689.9
$ws.Range("I80").Value = 803.125
$ws.Range("K80").Value = 2409.375
$ws.Range("M80").Value = -1411.375
$ws.Range("H83").Value = 3689.9
$ws.Range("I83").Value = 803.125
$ws.Range("K83").Value = 7228.125
$ws.Range("M83").Value = -2236.125
$ws.Range("H98").Value = 1872.375
$ws.Range("I98").Value = 1830.5333
$ws.Range("K98").Value = 1830.5333
$ws.Range("M98").Value = -332.5333000000001
$ws.Range("H122").Value = 1872.375
$ws.Range("I122").Value = 1830.5333
$ws.Range("K122").Value = 5491.5999
$ws.Range("M122").Value = -3041.5999
$ws.Range("H137").Value = 38791.27
$ws.Range("I137").Value = 52439
$ws.Range("K137").Value = 157317
$ws.Range("M137").Value = -154767
$ws.Range("H138").Value = 2918.3796
$ws.Range("I138").Value = 1805.5555
$ws.Range("J138").Value = 3496.1924
$ws.Range("K138").Value = 5416.666499999999
$ws.Range("L138").Value = 10488.5772
$ws.Range("M138").Value = -276.6664999999994
$ws.Range("N138").Value = -20768.5772
$ws.Range("H141").Value = 6736.8
$ws.Range("I141").Value = 7299.222
$ws.Range("K141").Value = 21897.666
$ws.Range("M141").Value = -16717.666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1586.25
$ws.Range("I32").Value = 1531.5657
$ws.Range("J32").Value = 7000
$ws.Range("K32").Value = 1531.5657
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = -1244.5657
$ws.Range("N32").Value = -7574
$ws.Range("H45").Value = 5105591.5
$ws.Range("I45").Value = 9525612
$ws.Range("K45").Value = 9525612
$ws.Range("M45").Value = -9525235

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 495.65216
$ws.Range("J80").Value = 463.91666
$ws.Range("L80").Value = 463.91666
$ws.Range("N80").Value = -2459.91666
$ws.Range("H83").Value = 495.65216
$ws.Range("J83").Value = 463.91666
$ws.Range("L83").Value = 2319.5833
$ws.Range("N83").Value = -12303.5833
$ws.Range("H86").Value = 33345848
$ws.Range("I86").Value = 108337530
$ws.Range("J86").Value = 16210.556
$ws.Range("K86").Value = 108337530
$ws.Range("L86").Value = 16210.556
$ws.Range("M86").Value = -108336407
$ws.Range("N86").Value = -18456.556
$ws.Range("H89").Value = 33345848
$ws.Range("I89").Value = 108337530
$ws.Range("J89").Value = 16210.556
$ws.Range("K89").Value = 541687650
$ws.Range("L89").Value = 81052.78
$ws.Range("M89").Value = -541682034
$ws.Range("N89").Value = -92284.78
$ws.Range("H94").Value = 2275602.8
$ws.Range("I94").Value = 2632617.8
$ws.Range("J94").Value = 14508.833
$ws.Range("K94").Value = 2632617.8
$ws.Range("L94").Value = 14508.833
$ws.Range("M94").Value = -2632166.8
$ws.Range("N94").Value = -15410.833
$ws.Range("H130").Value = 59750
$ws.Range("J130").Value = 59750
$ws.Range("L130").Value = 59750
$ws.Range("N130").Value = -69790
$ws.Range("H135").Value = 43519.668
$ws.Range("J135").Value = 43519.668
$ws.Range("L135").Value = 43519.668
$ws.Range("N135").Value = -53659.668
$ws.Range("H138").Value = 94363.75
$ws.Range("J138").Value = 94363.75
$ws.Range("L138").Value = 94363.75
$ws.Range("N138").Value = -104643.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29005.268
$ws.Range("I31").Value = 17643.5
$ws.Range("J31").Value = 30953
$ws.Range("K31").Value = 17643.5
$ws.Range("L31").Value = 30953
$ws.Range("M31").Value = -17348.5
$ws.Range("N31").Value = -31543
$ws.Range("H34").Value = 29005.268
$ws.Range("I34").Value = 17643.5
$ws.Range("J34").Value = 30953
$ws.Range("K34").Value = 17643.5
$ws.Range("L34").Value = 30953
$ws.Range("M34").Value = -17441.5
$ws.Range("N34").Value = -31357
$ws.Range("H62").Value = 7373
$ws.Range("I62").Value = 7373
$ws.Range("K62").Value = 7373
$ws.Range("M62").Value = -6749
$ws.Range("H65").Value = 7373
$ws.Range("I65").Value = 7373
$ws.Range("K65").Value = 36865
$ws.Range("M65").Value = -33745
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 116480.89
$ws.Range("J141").Value = 123719.72
$ws.Range("L141").Value = 123719.72
$ws.Range("N141").Value = -134079.72

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 773.5
$ws.Range("J108").Value = 250
$ws.Range("L108").Value = 750
$ws.Range("N108").Value = -6510
$ws.Range("H109").Value = 1680.1765
$ws.Range("I109").Value = 1877.091
$ws.Range("K109").Value = 5631.272999999999
$ws.Range("M109").Value = -4591.272999999999
$ws.Range("H110").Value = 13948.4
$ws.Range("I110").Value = 1845.2
$ws.Range("K110").Value = 5535.6
$ws.Range("M110").Value = -1445.6
$ws.Range("H111").Value = 1466.6666
$ws.Range("I111").Value = 1466.6666
$ws.Range("K111").Value = 4399.9998
$ws.Range("M111").Value = -1332.9998
$ws.Range("H112").Value = 449.66666
$ws.Range("I112").Value = 449.66666
$ws.Range("K112").Value = 1348.99998
$ws.Range("M112").Value = -240.9999800000001
$ws.Range("H133").Value = 2943.625
$ws.Range("I133").Value = 2943.625
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 8830.875
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -3770.875
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 1811.6428
$ws.Range("I134").Value = 1811.6428
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5434.928400000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -364.9284000000007
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 4502.933
$ws.Range("I137").Value = 2194.5715
$ws.Range("J137").Value = 6522.75
$ws.Range("K137").Value = 6583.7145
$ws.Range("L137").Value = 19568.25
$ws.Range("M137").Value = -1483.7145
$ws.Range("N137").Value = -29768.25
$ws.Range("H139").Value = 2626.6843
$ws.Range("I139").Value = 2305.7
$ws.Range("K139").Value = 6917.099999999999
$ws.Range("M139").Value = -1777.099999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1809.8182
$ws.Range("I68").Value = 741
$ws.Range("J68").Value = 2047.3334
$ws.Range("K68").Value = 741
$ws.Range("L68").Value = 2047.3334
$ws.Range("M68").Value = 8
$ws.Range("N68").Value = -3545.3334
$ws.Range("H71").Value = 1809.8182
$ws.Range("I71").Value = 741
$ws.Range("J71").Value = 2047.3334
$ws.Range("K71").Value = 3705
$ws.Range("L71").Value = 10236.667
$ws.Range("M71").Value = 39
$ws.Range("N71").Value = -17724.667
$ws.Range("H82").Value = 2527631.2
$ws.Range("I82").Value = 4632675.5
$ws.Range("J82").Value = 1578
$ws.Range("K82").Value = 4632675.5
$ws.Range("L82").Value = 1578
$ws.Range("M82").Value = -4632314.5
$ws.Range("N82").Value = -2300
$ws.Range("H85").Value = 2527631.2
$ws.Range("I85").Value = 4632675.5
$ws.Range("J85").Value = 1578
$ws.Range("K85").Value = 4632675.5
$ws.Range("L85").Value = 1578
$ws.Range("M85").Value = -4631427.5
$ws.Range("N85").Value = -4074
$ws.Range("H93").Value = 7755681.5
$ws.Range("I93").Value = 11906170
$ws.Range("J93").Value = 8103.8
$ws.Range("K93").Value = 11906170
$ws.Range("L93").Value = 8103.8
$ws.Range("M93").Value = -11904922
$ws.Range("N93").Value = -10599.8
$ws.Range("H136").Value = 124276.18
$ws.Range("I136").Value = 337616
$ws.Range("J136").Value = 7909
$ws.Range("K136").Value = 1012848
$ws.Range("L136").Value = 23727
$ws.Range("M136").Value = -1010298
$ws.Range("N136").Value = -28827

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11596.827
$ws.Range("I62").Value = 36313
$ws.Range("J62").Value = 8744.962
$ws.Range("K62").Value = 36313
$ws.Range("L62").Value = 8744.962
$ws.Range("M62").Value = -35689
$ws.Range("N62").Value = -9992.962
$ws.Range("H65").Value = 11596.827
$ws.Range("I65").Value = 36313
$ws.Range("J65").Value = 8744.962
$ws.Range("K65").Value = 181565
$ws.Range("L65").Value = 43724.81
$ws.Range("M65").Value = -178445
$ws.Range("N65").Value = -49964.81
$ws.Range("H107").Value = 34483764
$ws.Range("I107").Value = 58824236
$ws.Range("J107").Value = 1433.8334
$ws.Range("K107").Value = 176472708
$ws.Range("L107").Value = 4301.5002
$ws.Range("M107").Value = -176470788
$ws.Range("N107").Value = -8141.5002
